$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-17 18:22:21"

# Step 1: insert a fresh row at row 2 (pushes old rows 2-8 down to 3-9)
$ws.Rows("2:2").Insert()

# Step 2: insert two fresh rows at row 9 (pushes the old row 8, now at row 9, down to row 11)
$ws.Rows("9:10").Insert()

# Step 3: write the brand-new row 2 (new top entry)
$ws.Cells.Item(2,1).Value = $newTimestamp
$ws.Cells.Item(2,2).Value = "【急募】専門分野のAIチャットボットコードチェック依頼"
$ws.Cells.Item(2,3).Value = "システム開発"
$ws.Cells.Item(2,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(2,5).Value = "期限情報なし"
$ws.Cells.Item(2,6).Value = "https://www.lancers.jp/work/detail/5415270"
$ws.Cells.Item(2,7).Value = 303
$ws.Cells.Item(2,8).Value = "🔥AI,Ai"

# Step 4: write the two brand-new rows 9 and 10
$ws.Cells.Item(9,1).Value = $newTimestamp
$ws.Cells.Item(9,2).Value = "【クリエイティブ】Aurora Creative Lab 外注パートナー募集"
$ws.Cells.Item(9,3).Value = "システム開発"
$ws.Cells.Item(9,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(9,5).Value = "期限情報なし"
$ws.Cells.Item(9,6).Value = "https://www.lancers.jp/work/detail/5415615"
$ws.Cells.Item(9,7).Value = 18

$ws.Cells.Item(10,1).Value = $newTimestamp
$ws.Cells.Item(10,2).Value = "【再募集・急募・即決します】VBAで1問1答問題集の作成"
$ws.Cells.Item(10,3).Value = "システム開発"
$ws.Cells.Item(10,4).Value = "~ 5,000 円 / 固定"
$ws.Cells.Item(10,5).Value = "期限情報なし"
$ws.Cells.Item(10,6).Value = "https://www.lancers.jp/work/detail/5415538"
$ws.Cells.Item(10,7).Value = 10

# Step 5: refresh the timestamp (column A) on every existing data row (now rows 3-8 and 11)
$ws.Cells.Item(3,1).Value = $newTimestamp
$ws.Cells.Item(4,1).Value = $newTimestamp
$ws.Cells.Item(5,1).Value = $newTimestamp
$ws.Cells.Item(6,1).Value = $newTimestamp
$ws.Cells.Item(7,1).Value = $newTimestamp
$ws.Cells.Item(8,1).Value = $newTimestamp
$ws.Cells.Item(11,1).Value = $newTimestamp

# Step 6: the insert operations leave the worksheet's hyperlink list stale/misaligned,
# so rebuild hyperlinks for every URL cell from scratch.
$ws.Range("F2").Hyperlinks.Delete()

$urls = @{
    2 = "https://www.lancers.jp/work/detail/5415270"
    3 = "https://www.lancers.jp/work/detail/5415330"
    4 = "https://www.lancers.jp/work/detail/5415235"
    5 = "https://www.lancers.jp/work/detail/5405023"
    6 = "https://www.lancers.jp/work/detail/5398112"
    7 = "https://www.lancers.jp/work/detail/5415325"
    8 = "https://www.lancers.jp/work/detail/5415061"
    9 = "https://www.lancers.jp/work/detail/5415615"
    10 = "https://www.lancers.jp/work/detail/5415538"
    11 = "https://www.lancers.jp/work/detail/5414812"
}

foreach ($r in 2..11) {
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $urls[$r])
    $cell.Style = "Hyperlink"
}
